$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E100").Value = "Deep Learning for protein subcellular localisation prediction"
$ws.Range("E101").Value = "Narrative Emotion"
$ws.Range("E102").Value = "Deep Colourisation"
$ws.Range("E103").Value = "Data Analytics for AHA Case History Forms"
$ws.Range("E104").Value = "Forecasting Electricity Load for Commercial Buildings "
$ws.Range("E105").Value = "Standalone Software Metrics Tool"
$ws.Range("E106").Value = "Investigating Wavelet-Based Symbolic Representations for Time Series Classification "
$ws.Range("E107").Value = "Simulation of Hadoop Task Scheduling Algorithms in Distributed Computing Environments"
